# =====================================================================
# Adds player-info / extended-batting-stat scraping results as two new
# worksheets ("Player Info" and "ODI Batting Extra"), and reworks the
# MATCH_CARD_LINK column of the existing "ODI Batting" sheet into a
# plain MATCH_CODE column.
# =====================================================================

$wb = $excel.ActiveWorkbook
$odiBatting = $wb.Worksheets.Item("ODI Batting")

# ---------------------------------------------------------------------
# 1. Create the new sheets in the right order:
#    Player Info | ODI Batting | ODI Batting Extra
# ---------------------------------------------------------------------
$extra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $odiBatting)
$extra.Name = "ODI Batting Extra"

$info = $wb.Worksheets.Add($odiBatting)
$info.Name = "Player Info"

# NOTE: in this COM-interop implementation, worksheet variables behave
# like live references to a *position* in the Worksheets collection
# rather than a stable handle to a specific sheet object. Since sheets
# were just inserted/shuffled, every previously-captured handle
# ($odiBatting, $extra, ...) may now silently refer to the wrong sheet.
# Re-resolve every handle fresh, by name, now that the final sheet
# order/structure is settled, before doing any further editing.
$info = $wb.Worksheets.Item("Player Info")
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$extra = $wb.Worksheets.Item("ODI Batting Extra")

# ---------------------------------------------------------------------
# 2. Populate "Player Info"
# ---------------------------------------------------------------------
$infoHeader = $info.Range("A1:D1")
$infoHeader.Font.Bold = $true
$infoHeader.Borders.LineStyle = 1
$infoHeader.HorizontalAlignment = -4108
$infoHeader.VerticalAlignment = -4160

$info.Range("A1").Value = "ID"
$info.Range("B1").Value = "NAME"
$info.Range("C1").Value = "BATTING_HAND"
$info.Range("D1").Value = "BOWL_STYLE"

# ID looks numeric, so force text storage like the source data does
$info.Range("A2").NumberFormat = "@"
$info.Range("A2").Value = "4681"
$info.Range("B2").Value = "Shimron Odilon Hetmeyer"
$info.Range("C2").Value = "Left Handed"
$info.Range("D2").Value = "Left Arm Wrist Spin (Chinaman)"

# ---------------------------------------------------------------------
# 3. Rework "ODI Batting"
#    - D1 header: MATCH_CARD_LINK -> MATCH_CODE
#    - D2:D48: full scorecard URL -> bare match code
#    - Clear the stray empty B24 / B29 / B36 cells
# ---------------------------------------------------------------------
$odiBatting.Range("D1").Value = "MATCH_CODE"

$matchCodes = @{
    2  = "4100"
    3  = "4101"
    4  = "4144"
    5  = "4148"
    6  = "4150"
    7  = "4154"
    8  = "4159"
    9  = "4161"
    10 = "4164"
    11 = "4179"
    12 = "4180"
    13 = "4181"
    14 = "4213"
    15 = "4216"
    16 = "4219"
    17 = "4220"
    18 = "4221"
    19 = "4228"
    20 = "4229"
    21 = "4230"
    22 = "4253"
    23 = "4254"
    24 = "4255"
    25 = "4256"
    26 = "4260"
    27 = "4304"
    28 = "4312"
    29 = "4317"
    30 = "4321"
    31 = "4325"
    32 = "4333"
    33 = "4338"
    34 = "4344"
    35 = "4348"
    36 = "4359"
    37 = "4360"
    38 = "4362"
    39 = "4377"
    40 = "4378"
    41 = "4379"
    42 = "4385"
    43 = "4387"
    44 = "4388"
    45 = "4391"
    46 = "4394"
    47 = "4483"
    48 = "4486"
}

$codeColumn = $odiBatting.Range("D2:D48")
$codeColumn.NumberFormat = "@"
foreach ($row in 2..48) {
    $odiBatting.Range("D$row").Value = $matchCodes[$row]
}

$odiBatting.Range("B24").Value = ""
$odiBatting.Range("B29").Value = ""
$odiBatting.Range("B36").Value = ""

# ---------------------------------------------------------------------
# 4. Populate "ODI Batting Extra"
# ---------------------------------------------------------------------
$extraHeader = $extra.Range("A1:F1")
$extraHeader.Font.Bold = $true
$extraHeader.Borders.LineStyle = 1
$extraHeader.HorizontalAlignment = -4108
$extraHeader.VerticalAlignment = -4160

$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

# MATCH_CODE column looks numeric -> keep it textual like the source data
$extra.Range("A2:A21").NumberFormat = "@"
# NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL are stored as text in the source too
$extra.Range("C2:C21").NumberFormat = "@"
$extra.Range("D2:D21").NumberFormat = "@"
$extra.Range("E2:E21").NumberFormat = "@"

$extraRows = @(
    @{ Row = 2;  Code = "4317"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;     Mom = "NO" }
    @{ Row = 3;  Code = "4321"; Pos = 5;     N4 = "4";   N6 = "0";   Pct = "18.40%";  Mom = "NO" }
    @{ Row = 4;  Code = "4325"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;     Mom = "NO" }
    @{ Row = 5;  Code = "4333"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;     Mom = "NO" }
    @{ Row = 6;  Code = "4338"; Pos = 5;     N4 = "1";   N6 = "0";   Pct = "12.59%";  Mom = "NO" }
    @{ Row = 7;  Code = "4344"; Pos = 4;     N4 = "2";   N6 = "0";   Pct = "9.21%";   Mom = "NO" }
    @{ Row = 8;  Code = "4348"; Pos = 4;     N4 = "3";   N6 = "2";   Pct = "12.54%";  Mom = "NO" }
    @{ Row = 9;  Code = "4359"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;     Mom = "NO" }
    @{ Row = 10; Code = "4360"; Pos = 4;     N4 = "2";   N6 = "0";   Pct = "8.57%";   Mom = "NO" }
    @{ Row = 11; Code = "4362"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;     Mom = "NO" }
    @{ Row = 12; Code = "4377"; Pos = 3;     N4 = "0";   N6 = "0";   Pct = "1.52%";   Mom = "NO" }
    @{ Row = 13; Code = "4378"; Pos = 3;     N4 = "2";   N6 = "1";   Pct = "13.77%";  Mom = "NO" }
    @{ Row = 14; Code = "4379"; Pos = 3;     N4 = "0";   N6 = "0";   Pct = $null;     Mom = "NO" }
    @{ Row = 15; Code = "4385"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;     Mom = "NO" }
    @{ Row = 16; Code = "4387"; Pos = 3;     N4 = "0";   N6 = "0";   Pct = "1.43%";   Mom = "NO" }
    @{ Row = 17; Code = "4388"; Pos = 4;     N4 = "2";   N6 = "2";   Pct = "11.75%";  Mom = "NO" }
    @{ Row = 18; Code = "4391"; Pos = 4;     N4 = "2";   N6 = "0";   Pct = "4.35%";   Mom = "NO" }
    @{ Row = 19; Code = "4394"; Pos = 3;     N4 = "1";   N6 = "0";   Pct = "2.48%";   Mom = "NO" }
    @{ Row = 20; Code = "4483"; Pos = 2;     N4 = "2";   N6 = "0";   Pct = "8.94%";   Mom = "NO" }
    @{ Row = 21; Code = "4486"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;     Mom = "NO" }
)

foreach ($r in $extraRows) {
    $row = $r.Row
    $extra.Range("A$row").Value = $r.Code
    if ($null -ne $r.Pos) {
        $extra.Range("B$row").Value = $r.Pos
    }
    if ($null -ne $r.N4) {
        $extra.Range("C$row").Value = $r.N4
    }
    if ($null -ne $r.N6) {
        $extra.Range("D$row").Value = $r.N6
    }
    if ($null -ne $r.Pct) {
        $extra.Range("E$row").Value = $r.Pct
    }
    $extra.Range("F$row").Value = $r.Mom
}
